# Adds 4 new field rows (suspension_enriched, suspension_enriched_target,
# suspension_entity, suspension_entity_number) to the field-schema sheet,
# in their correct alphabetically-sorted position (row 268, just before
# "tissue_id"), each flagged for the "sample-suspension" assay (column AA),
# together with their explanatory cell comments. All fields that were
# previously on/after row 268 shift down by 4 rows to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkmark = [string][char]0x2713

# --- 1. Make room: insert 4 blank rows starting at row 268 -----------------
$ws.Rows("268:271").Insert()

# Row-insertion shifts cell VALUES down automatically, but cell COMMENTS in
# this runtime stay anchored to their original row number. So we must move
# the comments that belonged to the old rows 268-291 down to their new home
# at rows 272-295 ourselves (walking bottom-up so we never overwrite a
# comment before we've read it).
for ($r = 291; $r -ge 268; $r--) {
    $oldComment = $ws.Range("A$r").Comment
    $text = $oldComment.Text()
    $newRow = $r + 4
    $target = $ws.Range("A$newRow")
    if ($target.Comment -ne $null) {
        [void]$target.Comment.Text($text)
    } else {
        [void]$target.AddComment($text)
    }
}

# --- 2. Fill in the 4 new rows with field name + assay checkmark -----------
$ws.Range("A268").Value = "suspension_enriched"
$ws.Range("AA268").Value = $checkmark

$ws.Range("A269").Value = "suspension_enriched_target"
$ws.Range("AA269").Value = $checkmark

$ws.Range("A270").Value = "suspension_entity"
$ws.Range("AA270").Value = $checkmark

$ws.Range("A271").Value = "suspension_entity_number"
$ws.Range("AA271").Value = $checkmark

# --- 3. Set/replace the comments for the new rows ---------------------------
function Set-FieldComment($rowAddr, $text) {
    $rng = $ws.Range($rowAddr)
    if ($rng.Comment -ne $null) {
        [void]$rng.Comment.Text($text)
    } else {
        [void]$rng.AddComment($text)
    }
}

Set-FieldComment "A268" "Was the cell/nuclei population enriched?"
Set-FieldComment "A269" "If the suspension was enriched, then this is the target of the enrichment."
Set-FieldComment "A270" "The type of single cell entity derived from isolation protocol."
Set-FieldComment "A271" "Total number of cell/nuclei yielded post dissociation and enrichment."

Write-Host "Done. New dimension:" $ws.UsedRange.Rows.Count "rows"
